$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 221; this shifts rows 221:275 down to 222:276
$ws.Rows.Item(221).Insert()

# Copy row 222 (the row that used to be 221, now shifted down) formatting/values into new row 221
# then overwrite the cells that differ per the target data.
$ws.Range("A222:T222").Copy()
$ws.Range("A221").PasteSpecial()

$ws.Cells.Item(221, 4).Value = 44511        # D221 Fecha
$ws.Cells.Item(221, 11).Value = "Sin especificar"  # K221 Variedad
$ws.Cells.Item(221, 12).Value = "Primera Pintón"   # L221 Calidad
$ws.Cells.Item(221, 13).Value = 800          # M221 Volumen
$ws.Cells.Item(221, 14).Value = 22500        # N221 Precio minimo
$ws.Cells.Item(221, 15).Value = 23000        # O221 Precio maximo
$ws.Cells.Item(221, 16).Value = 22750        # P221 Precio promedio ponderado
$ws.Cells.Item(221, 19).Value = 1138         # S221 Precio $/Kg
